$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("N11").Value = 14758.31
$ws.Range("O11").Value = 14758.31

$ws.Range("N13").Value = 608625.34

$ws.Range("O14").Value = 221318.08

$ws.Range("O15").Value = 26968.24

$ws.Range("N20").Value = 39147.83

$ws.Range("K28").Value = 76128.56

$ws.Range("N30").Value = 2991.3

$wb.Save()
